# Insert a new data row at row 103 (weekly price update), shifting the
# existing rows 103:124 down to 104:125.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

$ws.Cells.Item(103, 1).Value = 8
$ws.Cells.Item(103, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(103, 3).Value = "Coquimbo"
$ws.Cells.Item(103, 4).Value = 44841
$ws.Cells.Item(103, 5).Value = 4
$ws.Cells.Item(103, 6).Value = 100112052
$ws.Cells.Item(103, 7).Value = "Albahaca"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 1200
$ws.Cells.Item(103, 11).Value = 4000
$ws.Cells.Item(103, 12).Value = 4500
$ws.Cells.Item(103, 13).Value = 4250
$ws.Cells.Item(103, 14).Value = "`$/paquete"
$ws.Cells.Item(103, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(103, 16).Value = 4250
$ws.Cells.Item(103, 17).Value = 1
$ws.Cells.Item(103, 18).Value = "Hortaliza"
